# Revert 'cards' to commit 95cda46ab8 (Jun 25)
# account_type.xlsx: bring back the Airbyte sync metadata columns around the
# existing (aty_code, aty_labe, updated_at) data, matching the pre-revert
# layout: _airbyte_ab_id, _airbyte_emitted_at, aty_code, aty_labe,
# _airbyte_additional_properties, source_file_path, updated_at.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two leading Airbyte columns -----------------------------
# Before: A=aty_code, B=aty_labe, C=updated_at
$ws.Range("A:B").Insert()
# After:  A=(new) B=(new) C=aty_code D=aty_labe E=updated_at

# --- Insert the two trailing Airbyte columns (between aty_labe & updated_at) ---
$ws.Range("E:F").Insert()
# After:  A=(new) B=(new) C=aty_code D=aty_labe E=(new) F=(new) G=updated_at

# --- Header row values -----------------------------------------------------
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("E1").Value = "_airbyte_additional_properties"
$ws.Range("F1").Value = "source_file_path"

# The E:F insert already inherited the bold/border/centered header format
# from its neighbours, but the A:B insert (leftmost column, no left
# neighbour) did not - so explicitly clone the header format onto A1:B1
# from the (now) "aty_code" header cell.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data rows ---------------------------------------------------------
$abIds = @(
    "ad969c63-4ef2-45c3-847b-85f002584c88",
    "4a38db5d-dc2e-4eea-888a-78c6c3b94f6a",
    "cae9d427-1b55-4414-9544-9f50c1d392c1",
    "959f9db3-5a06-4bb9-8b88-96498171a321"
)
$sourceFiles = @(
    "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/ACCOUNT_TYPE/2024_08_06_1722929004063_1.parquet",
    "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/ACCOUNT_TYPE/2024_08_06_1722929004063_1.parquet",
    "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/ACCOUNT_TYPE/2024_08_06_1722929004063_1.parquet",
    "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/ACCOUNT_TYPE/2024_08_06_1722929004063_0.parquet"
)

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2

    $ws.Cells.Item($r, 1).Value = $abIds[$i]
    $ws.Cells.Item($r, 2).Value = 45510.3079196875
    $ws.Cells.Item($r, 5).Value = ""
    $ws.Cells.Item($r, 6).Value = $sourceFiles[$i]
}

# Give the new _airbyte_emitted_at column (B) the same date/time number
# format already used by updated_at (G), so it matches style-for-style.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("B2:B5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A1").Select()
